# Generate Report for Handoff
#
# Swaps which localization job is shown first on each sheet: the
# "c000875b-...md" file moves into row 2 (still "In Translation") and the
# "90e98aa2-...md" file moves into row 3, whose status flips to
# "Ready for handoff" now that a fresh handoff has just been generated
# (new Latest Handoff Datetime stamps on the zh-cn / de-de detail sheets).

$wb = $excel.ActiveWorkbook

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/a798acfc28fe372c173a0cecb9e887509bb56135"
$urlC000875bMd = "$mdBase/e2e/c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md"
$urlC000875bMd = "https://github.com/OpenLocalizationTest/oltest/blob/a798acfc28fe372c173a0cecb9e887509bb56135/e2e/c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md"
$url90e98aa2Md = "https://github.com/OpenLocalizationTest/oltest/blob/a798acfc28fe372c173a0cecb9e887509bb56135/e2e/90e98aa2-61c1-45b2-a382-4df0ea003a41.md"
$urlLocConfig  = "https://github.com/OpenLocalizationTest/oltest/blob/a798acfc28fe372c173a0cecb9e887509bb56135/.localization-config"

$urlC000875bZhCnXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16c59acf7bb853578ebb982c306d8e2eb69edf65/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/c000875b-8a7c-45c6-9bce-2eebeb2b71ee.3013e880563e2b4a6fbb986dd1a6855a0c29e4b0.zh-cn.xlf"
$url90e98aa2ZhCnXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/16c59acf7bb853578ebb982c306d8e2eb69edf65/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/90e98aa2-61c1-45b2-a382-4df0ea003a41.0549efd166e76ef3e2d649eb8d4dea2a8a6339ec.zh-cn.xlf"

$urlC000875bDeDeXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1851a1798be2064f0c6de277ad1957d2cb632e74/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/c000875b-8a7c-45c6-9bce-2eebeb2b71ee.3013e880563e2b4a6fbb986dd1a6855a0c29e4b0.de-de.xlf"
$url90e98aa2DeDeXlf = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/1851a1798be2064f0c6de277ad1957d2cb632e74/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/90e98aa2-61c1-45b2-a382-4df0ea003a41.0549efd166e76ef3e2d649eb8d4dea2a8a6339ec.de-de.xlf"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Range.Hyperlinks.Delete() clears every hyperlink on the sheet, so do it
# once up front and rebuild all of them afterwards.
$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "In Translation"

$ws.Range("A3").Value = "90e98aa2-61c1-45b2-a382-4df0ea003a41.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("C4").Value = "Not to be localized"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlC000875bMd, "", "", "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md")
$ws.Hyperlinks.Add($ws.Range("A3"), $url90e98aa2Md, "", "", "90e98aa2-61c1-45b2-a382-4df0ea003a41.md")
$ws.Hyperlinks.Add($ws.Range("A4"), $urlLocConfig, "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.3013e880563e2b4a6fbb986dd1a6855a0c29e4b0.zh-cn.xlf"
$ws.Range("D2").Value = "2016-03-08 08:11:26"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "90e98aa2-61c1-45b2-a382-4df0ea003a41.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "90e98aa2-61c1-45b2-a382-4df0ea003a41.0549efd166e76ef3e2d649eb8d4dea2a8a6339ec.zh-cn.xlf"
$ws.Range("D3").Value = "2016-03-08 08:12:04"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlC000875bMd, "", "", "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $urlC000875bZhCnXlf, "", "", "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.3013e880563e2b4a6fbb986dd1a6855a0c29e4b0.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $url90e98aa2Md, "", "", "90e98aa2-61c1-45b2-a382-4df0ea003a41.md")
$ws.Hyperlinks.Add($ws.Range("C3"), $url90e98aa2ZhCnXlf, "", "", "90e98aa2-61c1-45b2-a382-4df0ea003a41.0549efd166e76ef3e2d649eb8d4dea2a8a6339ec.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), $urlLocConfig, "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Hyperlinks.Delete()

$ws.Range("A2").Value = "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md"
$ws.Range("B2").Value = "In Translation"
$ws.Range("C2").Value = "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.3013e880563e2b4a6fbb986dd1a6855a0c29e4b0.de-de.xlf"
$ws.Range("D2").Value = "2016-03-08 08:11:31"
$ws.Range("G2").Value = "0001-01-01 00:00:00"
$ws.Range("H2").Value = "Include"

$ws.Range("A3").Value = "90e98aa2-61c1-45b2-a382-4df0ea003a41.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "90e98aa2-61c1-45b2-a382-4df0ea003a41.0549efd166e76ef3e2d649eb8d4dea2a8a6339ec.de-de.xlf"
$ws.Range("D3").Value = "2016-03-08 08:12:08"
$ws.Range("G3").Value = "0001-01-01 00:00:00"
$ws.Range("H3").Value = "Include"

$ws.Range("A4").Value = ".localization-config"
$ws.Range("B4").Value = "Not to be localized"
$ws.Range("D4").Value = "0001-01-01 00:00:00"
$ws.Range("G4").Value = "0001-01-01 00:00:00"
$ws.Range("H4").Value = "Ignored"

$ws.Hyperlinks.Add($ws.Range("A2"), $urlC000875bMd, "", "", "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.md")
$ws.Hyperlinks.Add($ws.Range("C2"), $urlC000875bDeDeXlf, "", "", "c000875b-8a7c-45c6-9bce-2eebeb2b71ee.3013e880563e2b4a6fbb986dd1a6855a0c29e4b0.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), $url90e98aa2Md, "", "", "90e98aa2-61c1-45b2-a382-4df0ea003a41.md")
$ws.Hyperlinks.Add($ws.Range("C3"), $url90e98aa2DeDeXlf, "", "", "90e98aa2-61c1-45b2-a382-4df0ea003a41.0549efd166e76ef3e2d649eb8d4dea2a8a6339ec.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), $urlLocConfig, "", "", ".localization-config")
